# Rename the "Flux units" header on the Reactions sheet to "Flux bound units"
$wb = $excel.ActiveWorkbook

$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Range("H1").Value = "Flux bound units"
$wsReactions.Activate()
$wsReactions.Range("H1").Select()

# Insert two new columns on the "dFBA objectives" sheet for the new
# "Reaction rate units" / "Coefficient units" headers (reflecting the
# wc_lang change to reaction rate units), ahead of the existing
# "Database references" column.
$wsObjectives = $wb.Worksheets.Item("dFBA objectives")
$wsObjectives.Range("F1:G1").EntireColumn.Insert()
$wsObjectives.Cells.Item(1, 6).Value = "Reaction rate units"
$wsObjectives.Cells.Item(1, 7).Value = "Coefficient units"

# Leave this sheet active, as it was the last one the author touched.
$wsObjectives.Activate()
$wsObjectives.Range("I11").Select()
